# Update basaltic glass analysis - major, wd scans and D2872
# Applies corrected values to D2872 kraw summary sheet (rows 2-15)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 100.07
$ws.Range("I2").Value = 100.07
$ws.Range("J2").Value = 0.96
$ws.Range("K2").Value = 5.3
$ws.Range("M2").Value = 6.42
$ws.Range("E3").Value = 100.07
$ws.Range("H3").Value = 1.23
$ws.Range("I3").Value = 100.07
$ws.Range("J3").Value = 1.02
$ws.Range("K3").Value = 4.58
$ws.Range("L3").Value = 1.23
$ws.Range("M3").Value = 5.84
$ws.Range("E4").Value = 100.07
$ws.Range("H4").Value = 1.04
$ws.Range("I4").Value = 100.07
$ws.Range("K4").Value = 6.55
$ws.Range("L4").Value = 1.04
$ws.Range("M4").Value = 7.48
$ws.Range("E5").Value = 100.07
$ws.Range("H5").Value = 0.76
$ws.Range("I5").Value = 100.07
$ws.Range("J5").Value = 0.63
$ws.Range("K5").Value = 20.48
$ws.Range("L5").Value = 0.76
$ws.Range("M5").Value = 20.79
$ws.Range("E6").Value = 100.07
$ws.Range("I6").Value = 100.07
$ws.Range("J6").Value = 0.91
$ws.Range("K6").Value = 6.52
$ws.Range("M6").Value = 7.45
$ws.Range("E7").Value = 100.07
$ws.Range("I7").Value = 100.07
$ws.Range("J7").Value = 0.73
$ws.Range("K7").Value = 7.46
$ws.Range("M7").Value = 8.289999999999999
$ws.Range("E8").Value = 100.07
$ws.Range("I8").Value = 100.07
$ws.Range("J8").Value = 0.6899999999999999
$ws.Range("K8").Value = 9.26
$ws.Range("M8").Value = 9.94
$ws.Range("E9").Value = 100.07
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 100.07
$ws.Range("J9").Value = 0.83
$ws.Range("K9").Value = 7.29
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 8.140000000000001
$ws.Range("E10").Value = 100.07
$ws.Range("H10").Value = 0.89
$ws.Range("I10").Value = 100.07
$ws.Range("J10").Value = 0.74
$ws.Range("K10").Value = 6.97
$ws.Range("L10").Value = 0.89
$ws.Range("M10").Value = 7.85
$ws.Range("E11").Value = 100.07
$ws.Range("H11").Value = 0.9399999999999999
$ws.Range("I11").Value = 100.07
$ws.Range("J11").Value = 0.78
$ws.Range("K11").Value = 16.33
$ws.Range("L11").Value = 0.9399999999999999
$ws.Range("M11").Value = 16.72
$ws.Range("E12").Value = 100.07
$ws.Range("H12").Value = 1.05
$ws.Range("I12").Value = 100.07
$ws.Range("J12").Value = 0.87
$ws.Range("K12").Value = 15.55
$ws.Range("L12").Value = 1.05
$ws.Range("M12").Value = 15.96
$ws.Range("E13").Value = 100.07
$ws.Range("H13").Value = 0.8
$ws.Range("I13").Value = 100.07
$ws.Range("J13").Value = 0.66
$ws.Range("K13").Value = 19.36
$ws.Range("L13").Value = 0.8
$ws.Range("M13").Value = 19.69
$ws.Range("E14").Value = 100.07
$ws.Range("H14").Value = 1.04
$ws.Range("I14").Value = 100.07
$ws.Range("J14").Value = 0.86
$ws.Range("K14").Value = 12.71
$ws.Range("L14").Value = 1.04
$ws.Range("M14").Value = 13.21
$ws.Range("E15").Value = 100.07
$ws.Range("H15").Value = 1.04
$ws.Range("I15").Value = 100.07
$ws.Range("J15").Value = 0.86
$ws.Range("K15").Value = 14.23
$ws.Range("L15").Value = 1.04
$ws.Range("M15").Value = 14.68